# CONTROL_MAPPINGS.xlsx — fix control-signal formulas + register mapping table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column G (REG-DEST / ALU-SRC style bit): was
#   IF(AND(OR(Fn="R", Fn="I"),  NOT(En="lw")), 1, 0)
# now
#   IF(AND(OR(Fn="R", En="sw"), NOT(En="sll"),  NOT(En="srl")), 1, 0)
# applied for every instruction row (3-18)
# ---------------------------------------------------------------------------
for ($r = 3; $r -le 18; $r++) {
    $ws.Cells.Item($r, 7).Formula = '=IF(AND(OR(F' + $r + '="R", E' + $r + '="sw"), NOT(E' + $r + '="sll"),  NOT(E' + $r + '="srl")), 1, 0)'
}

# ---------------------------------------------------------------------------
# Column J: was (buggy - referenced row+2 instead of the current row)
#   IF(OR(F(n+2)="I", F(n+2)="R", E(n+2)="lw"), 1, 0)
# now
#   IF(AND(OR(Fn="I", Fn="R", En="lw"), NOT(En="sw"), NOT(En="beq"), NOT(En="bneq")), 1, 0)
# ---------------------------------------------------------------------------
for ($r = 3; $r -le 18; $r++) {
    $ws.Cells.Item($r, 10).Formula = '=IF(AND(OR(F' + $r + '="I", F' + $r + '="R", E' + $r + '="lw"), NOT(E' + $r + '="sw"), NOT(E' + $r + '="beq"), NOT(E' + $r + '="bneq")), 1, 0)'
}

# ---------------------------------------------------------------------------
# Column R (ALUOP bit): was
#   IF(OR(TRIM(En)="sub", TRIM(En)="subi", TRIM(En)="and", TRIM(En)="andi", TRIM(En)="sll", TRIM(En)="sra"), 1, 0)
# now also includes bneq/beq
#   IF(OR(TRIM(En)="sub", TRIM(En)="subi", TRIM(En)="and", TRIM(En)="andi", TRIM(En)="sll", TRIM(En)="sra", TRIM(En)="bneq", TRIM(En)="beq"), 1, 0)
# ---------------------------------------------------------------------------
for ($r = 3; $r -le 18; $r++) {
    $ws.Cells.Item($r, 18).Formula = '=IF(OR(TRIM(E' + $r + ')="sub", TRIM(E' + $r + ')="subi", TRIM(E' + $r + ')="and", TRIM(E' + $r + ')="andi", TRIM(E' + $r + ')="sll", TRIM(E' + $r + ')="sra", TRIM(E' + $r + ')="bneq", TRIM(E' + $r + ')="beq"), 1, 0)'
}

# ---------------------------------------------------------------------------
# Register-mapping side table (W6:W9): rotate $t4,$t3,$t2,$t1 -> $t1,$t2,$t3,$t4
# (single-quoted literals -- these contain a literal "$" which must not be
# treated as PowerShell variable interpolation)
# ---------------------------------------------------------------------------
$ws.Range("W6").Value = '$t1'
$ws.Range("W7").Value = '$t2'
$ws.Range("W8").Value = '$t3'
$ws.Range("W9").Value = '$t4'

# ---------------------------------------------------------------------------
# Recalculate so every dependent cell (S,T,U columns + the rolled-up U20) has
# a fresh cached value, then move the active selection to J18 (matches the
# final cursor position saved with the workbook).
# ---------------------------------------------------------------------------
$wb.Application.Calculate()
$ws.Range("J18").Select()
